$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 171, pushing the existing rows
# (old 171-179) down to 173-181.
$ws.Rows.Item(171).Insert()
$ws.Rows.Item(171).Insert()

# New row 171
$ws.Range("A171").Value = 1
$ws.Range("B171").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C171").Value = "Arica y Parinacota"
$ws.Range("D171").Value = 45008
$ws.Range("E171").Value = 15
$ws.Range("F171").Value = 100114001
$ws.Range("G171").Value = "Papa"
$ws.Range("H171").Value = "Asterix"
$ws.Range("I171").Value = "1a (cosecha)"
$ws.Range("J171").Value = 1000
$ws.Range("K171").Value = 13000
$ws.Range("L171").Value = 14000
$ws.Range("M171").Value = 13500
$ws.Range("N171").Value = "$/saco 25 kilos"
$ws.Range("O171").Value = "Región de Los Lagos"
$ws.Range("P171").Value = 540
$ws.Range("Q171").Value = 25
$ws.Range("R171").Value = "Hortaliza"

# New row 172
$ws.Range("A172").Value = 1
$ws.Range("B172").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C172").Value = "Arica y Parinacota"
$ws.Range("D172").Value = 45008
$ws.Range("E172").Value = 15
$ws.Range("F172").Value = 100114001
$ws.Range("G172").Value = "Papa"
$ws.Range("H172").Value = "Red Lady"
$ws.Range("I172").Value = "1a (cosecha)"
$ws.Range("J172").Value = 1000
$ws.Range("K172").Value = 12000
$ws.Range("L172").Value = 13000
$ws.Range("M172").Value = 12500
$ws.Range("N172").Value = "$/saco 25 kilos"
$ws.Range("O172").Value = "Región del Maule"
$ws.Range("P172").Value = 500
$ws.Range("Q172").Value = 25
$ws.Range("R172").Value = "Hortaliza"
